# Update cryptocurrency price (D) and volume change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format so price strings like "1.005"
# or "27.825.14" are not reinterpreted as numbers/dates by Excel
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.825.14"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "1.765.76"
$ws.Range("E3").Value = "  -2.62%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.69%  "
$ws.Range("D5").Value = "339.25"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("D7").Value = "0.3772"
$ws.Range("E7").Value = "  -3.94%  "
$ws.Range("D8").Value = "0.3358"
$ws.Range("E8").Value = "  -3.81%  "
$ws.Range("D9").Value = "45.61"
$ws.Range("E9").Value = "  -5.49%  "
$ws.Range("D10").Value = "1.134"
$ws.Range("E10").Value = "  -5.72%  "
$ws.Range("D11").Value = "0.07261"
$ws.Range("E11").Value = "  -4.26%  "
$ws.Range("D12").Value = "23.01"
$ws.Range("E12").Value = "  +3.75%  "
$ws.Range("D13").Value = "1.004"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").Value = "6.251"
$ws.Range("E14").Value = "  -4.26%  "
$ws.Range("D15").Value = "7.234"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").Value = "1.766.77"
$ws.Range("E16").Value = "  -2.51%  "
$ws.Range("D17").Value = "0.00001056"
$ws.Range("E17").Value = "  -4.75%  "
$ws.Range("D18").Value = "0.06605"
$ws.Range("E18").Value = "  -1.11%  "
$ws.Range("D19").Value = "81.07"
$ws.Range("E19").Value = "  -4.83%  "
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("D21").Value = "17.15"
$ws.Range("E21").Value = "  -4.18%  "
$ws.Range("D22").Value = "6.328"
$ws.Range("E22").Value = "  -3.85%  "
$ws.Range("D23").Value = "27.828.52"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").Value = "11.77"
$ws.Range("E24").Value = "  -8.62%  "
$ws.Range("D25").Value = "2.380"
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("D26").Value = "1.489"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").Value = "152.60"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("D28").Value = "19.99"
$ws.Range("E28").Value = "  -6.44%  "
$ws.Range("D29").Value = "2.358"
$ws.Range("E29").Value = "  -7.94%  "
$ws.Range("D30").Value = "1.969.12"
$ws.Range("E30").Value = "  -2.40%  "
$ws.Range("D31").Value = "132.27"
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("D32").Value = "4.035"
$ws.Range("D33").Value = "5.910"
$ws.Range("E33").Value = "  -4.05%  "
$ws.Range("D34").Value = "0.08751"
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("D35").Value = "12.42"
$ws.Range("E35").Value = "  -6.81%  "
$ws.Range("D36").Value = "0.02361"
$ws.Range("E36").Value = "  -2.94%  "
$ws.Range("D37").Value = "0.6709"
$ws.Range("E37").Value = "  -3.24%  "
$ws.Range("D38").Value = "0.06266"
$ws.Range("E38").Value = "  -4.39%  "
$ws.Range("D39").Value = "5.208"
$ws.Range("E39").Value = "  -6.32%  "
$ws.Range("D40").Value = "0.2119"
$ws.Range("E40").Value = "  -4.95%  "
$ws.Range("D41").Value = "1.224"
$ws.Range("E41").Value = "  -3.39%  "
$ws.Range("D42").Value = "1.474"
$ws.Range("E42").Value = "  -8.52%  "
$ws.Range("D43").Value = "8.089"
$ws.Range("E43").Value = "  -5.73%  "
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").Value = "13.92"
$ws.Range("E45").Value = "  -5.28%  "
$ws.Range("D46").Value = "0.6119"
$ws.Range("E46").Value = "  -6.71%  "
$ws.Range("D47").Value = "3.844"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("D48").Value = "131.45"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("D49").Value = "2.027"
$ws.Range("E49").Value = "  -6.46%  "
$ws.Range("D50").Value = "0.07277"
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("D51").Value = "1.183"
$ws.Range("E51").Value = "  +1.62%  "

# Restore the original (default) style on column D now that the text is set
$ws.Range("D2:D51").Style = "Normal"
